$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, pushing existing rows 51-61 down to 52-62.
$ws.Rows(51).Insert()

# Populate the newly inserted row 51 with the new weekly price record.
$ws.Cells.Item(51, 1).Value = 10
$ws.Cells.Item(51, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(51, 3).Value = "La Araucanía"
$ws.Cells.Item(51, 4).Value = 44522
$ws.Cells.Item(51, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(51, 5).Value = 9
$ws.Cells.Item(51, 6).Value = 300000000
$ws.Cells.Item(51, 7).Value = "Espárragos"
$ws.Cells.Item(51, 8).Value = "Sin especificar"
$ws.Cells.Item(51, 9).Value = "Primera"
$ws.Cells.Item(51, 10).Value = 400
$ws.Cells.Item(51, 11).Value = 1300
$ws.Cells.Item(51, 12).Value = 1400
$ws.Cells.Item(51, 13).Value = 1350
$ws.Cells.Item(51, 14).Value = "$/kilo"
$ws.Cells.Item(51, 15).Value = "Región del Maule"
$ws.Cells.Item(51, 16).Value = 1350
$ws.Cells.Item(51, 17).Value = 1
$ws.Cells.Item(51, 18).Value = "Hortaliza"
